# Versiebeheer.docx — add "Versie 1.2" and "Versie 1.1" changelog
# sections above the existing "Versie 1.0" section (commit: "Add more
# options to sabotaging the air conditioning").

$d = $word.ActiveDocument

# The green used for every changelog bullet elsewhere in the document
# (OOXML hex 00B050) expressed as the BGR-packed wdColor long Word uses.
$greenColor = 0x50B000

# --- 1. Make room: insert 7 empty paragraphs above the current
#        first paragraph ("Versie 1.0"), one per new line we need. ---
$firstPara = $d.Paragraphs(1).Range
for ($i = 0; $i -lt 7; $i++) {
    $firstPara.InsertParagraphBefore()
}

# After the inserts, paragraphs 1-7 are the new (still empty) ones and
# paragraph 8 is the original "Versie 1.0" heading.

# --- 2. Versie 1.2 heading ---
$p = $d.Paragraphs(1)
$p.Range.Text = "Versie 1.2"
$p.Style = "Kop1"

# --- 3. The four Versie 1.2 bullets (green, "Lijstalinea" list style) ---
$p = $d.Paragraphs(2)
$p.Range.Text = "Bugfixes"
$p.Style = "Lijstalinea"
$p.Range.Font.Color = $greenColor

$p = $d.Paragraphs(3)
$p.Range.Text = "Bevestiging toevoegen bij het spel verlaten"
$p.Style = "Lijstalinea"
$p.Range.Font.Color = $greenColor

$p = $d.Paragraphs(4)
$p.Range.Text = "Meer aandacht voor prestatie bereikt"
$p.Style = "Lijstalinea"
$p.Range.Font.Color = $greenColor

$p = $d.Paragraphs(5)
$p.Range.Text = "Meer opties toevoegen aan airco saboteren"
$p.Style = "Lijstalinea"
$p.Range.Font.Color = $greenColor

# --- 4. Versie 1.1 heading ---
$p = $d.Paragraphs(6)
$p.Range.Text = "Versie 1.1"
$p.Style = "Kop1"

# --- 5. The single Versie 1.1 bullet ---
$p = $d.Paragraphs(7)
$p.Range.Text = "Bugfixes"
$p.Style = "Lijstalinea"
$p.Range.Font.Color = $greenColor

# --- 6. Give the five new bullets one shared list numbering instance,
#        copied from the list style already used throughout the rest
#        of the document (same underlying "Lijststijlvormen" list). ---
$existingListTemplate = $d.Paragraphs(12).Range.ListFormat.ListTemplate
$bulletsRange = $d.Range($d.Paragraphs(2).Range.Start, $d.Paragraphs(5).Range.End)
$bulletsRange.ListFormat.ApplyListTemplate($existingListTemplate)
$lastBulletRange = $d.Paragraphs(7).Range
$lastBulletRange.ListFormat.ApplyListTemplate($existingListTemplate)

# --- 7. Move the document's "_GoBack" bookmark (Word always keeps at
#        most one) to sit inside the word "airco", matching the last
#        edit position recorded by the author. ---
$p5 = $d.Paragraphs(5).Range
$bmPos = $p5.Start + 30
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output ("Done. Paragraph count=" + $d.Paragraphs.Count)
